$wb = $excel.ActiveWorkbook

# Sheet names that contain the exhibition data table that needs updating
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 2212
    $ws.Range("F4").Value = 346
    $ws.Range("F6").Value = 6421
    $ws.Range("F7").Value = 290
}
